$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows (including the
# "Lakes Entrance" block) down by one.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new "Emerald" exposure site.
$ws.Cells.Item(8, 1).Value = "Emerald"
$ws.Cells.Item(8, 2).Value = "Lakeside Paddle Boats, Emerald Lake Park"
$ws.Cells.Item(8, 3).Value = "31/12/20 3:30pm - 5:30pm"
$ws.Cells.Item(8, 4).Value = "Case visited venue"

# The insert pushed the old last "Melbourne" row (Nandos entry) down into
# row 22, clashing with the following "Nunawading" row. Remove that
# displaced duplicate row so the rest of the table keeps its original
# row positions.
$ws.Rows.Item(22).Delete()
